# Weekly roll of the Rabanito / Vega Modelo de Temuco price log:
# a new week's record is inserted at row 18 (pushing the existing
# rows 18-61 down to 19-62), populated with the same market/category
# metadata as the (now shifted) row 19, but with this week's date and
# reported volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 18; this shifts old rows
# 18..61 down to 19..62 (and bumps the sheet's used-range dimension
# to R62 automatically).
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new week's record.
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44690
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 300000001
$ws.Range("G18").Value = "Rabanito"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("N18").Value = "$/docena de paquetes"
$ws.Range("O18").Value = "Provincia de Cautín"
$ws.Range("P18").Value = 583
$ws.Range("Q18").Value = 12
$ws.Range("R18").Value = "Hortaliza"
